$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-06-01 Saturday"

# Update every cell of the practice table with newly generated problems
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "21-4="
$t.Cell(1, 2).Range.Text = "44+29="
$t.Cell(1, 3).Range.Text = "90-36="
$t.Cell(1, 4).Range.Text = "92-28="
$t.Cell(1, 5).Range.Text = "40-9="

$t.Cell(2, 1).Range.Text = "65-16="
$t.Cell(2, 2).Range.Text = "9+55="
$t.Cell(2, 3).Range.Text = "27+26="
$t.Cell(2, 4).Range.Text = "58+15="
$t.Cell(2, 5).Range.Text = "22-17="

$t.Cell(3, 1).Range.Text = "82-38="
$t.Cell(3, 2).Range.Text = "7+34="
$t.Cell(3, 3).Range.Text = "40-4="
$t.Cell(3, 4).Range.Text = "64-57="
$t.Cell(3, 5).Range.Text = "17+7="

$t.Cell(4, 1).Range.Text = "66+5="
$t.Cell(4, 2).Range.Text = "27-18="
$t.Cell(4, 3).Range.Text = "95-79="
$t.Cell(4, 4).Range.Text = "17+15="
$t.Cell(4, 5).Range.Text = "47-9="

$t.Cell(5, 1).Range.Text = "18+77="
$t.Cell(5, 2).Range.Text = "25+8="
$t.Cell(5, 3).Range.Text = "60-46="
$t.Cell(5, 4).Range.Text = "91-24="
$t.Cell(5, 5).Range.Text = "8+64="

$t.Cell(6, 1).Range.Text = "91-52="
$t.Cell(6, 2).Range.Text = "38-29="
$t.Cell(6, 3).Range.Text = "68+6="
$t.Cell(6, 4).Range.Text = "91-33="
$t.Cell(6, 5).Range.Text = "15+78="

$t.Cell(7, 1).Range.Text = "41-24="
$t.Cell(7, 2).Range.Text = "50-37="
$t.Cell(7, 3).Range.Text = "82-24="
$t.Cell(7, 4).Range.Text = "9+46="
$t.Cell(7, 5).Range.Text = "58+18="

$t.Cell(8, 1).Range.Text = "37+19="
$t.Cell(8, 2).Range.Text = "81-48="
$t.Cell(8, 3).Range.Text = "34+9="
$t.Cell(8, 4).Range.Text = "15+68="
$t.Cell(8, 5).Range.Text = "45+47="

$t.Cell(9, 1).Range.Text = "26+25="
$t.Cell(9, 2).Range.Text = "65+27="
$t.Cell(9, 3).Range.Text = "39+6="
$t.Cell(9, 4).Range.Text = "19+4="
$t.Cell(9, 5).Range.Text = "85-37="

$t.Cell(10, 1).Range.Text = "47+14="
$t.Cell(10, 2).Range.Text = "51-28="
$t.Cell(10, 3).Range.Text = "47+26="
$t.Cell(10, 4).Range.Text = "50-23="
$t.Cell(10, 5).Range.Text = "71-27="

$t.Cell(11, 1).Range.Text = "68+8="
$t.Cell(11, 2).Range.Text = "5+16="
$t.Cell(11, 3).Range.Text = "47+6="
$t.Cell(11, 4).Range.Text = "71-23="
$t.Cell(11, 5).Range.Text = "80-42="

$t.Cell(12, 1).Range.Text = "56+39="
$t.Cell(12, 2).Range.Text = "74-7="
$t.Cell(12, 3).Range.Text = "44-36="
$t.Cell(12, 4).Range.Text = "44-5="
$t.Cell(12, 5).Range.Text = "48+25="

$t.Cell(13, 1).Range.Text = "33-4="
$t.Cell(13, 2).Range.Text = "9+48="
$t.Cell(13, 3).Range.Text = "56-29="
$t.Cell(13, 4).Range.Text = "64-58="
$t.Cell(13, 5).Range.Text = "42-34="

$t.Cell(14, 1).Range.Text = "92-35="
$t.Cell(14, 2).Range.Text = "18+17="
$t.Cell(14, 3).Range.Text = "78+3="
$t.Cell(14, 4).Range.Text = "92-8="
$t.Cell(14, 5).Range.Text = "60-18="

$t.Cell(15, 1).Range.Text = "95-48="
$t.Cell(15, 2).Range.Text = "70-3="
$t.Cell(15, 3).Range.Text = "9+73="
$t.Cell(15, 4).Range.Text = "7+59="
$t.Cell(15, 5).Range.Text = "23+58="

$t.Cell(16, 1).Range.Text = "91-35="
$t.Cell(16, 2).Range.Text = "38+14="
$t.Cell(16, 3).Range.Text = "97-78="
$t.Cell(16, 4).Range.Text = "91-87="
$t.Cell(16, 5).Range.Text = "26+46="

$t.Cell(17, 1).Range.Text = "33-6="
$t.Cell(17, 2).Range.Text = "48+9="
$t.Cell(17, 3).Range.Text = "30-8="
$t.Cell(17, 4).Range.Text = "70-39="
$t.Cell(17, 5).Range.Text = "97-19="

$t.Cell(18, 1).Range.Text = "67-19="
$t.Cell(18, 2).Range.Text = "90-48="
$t.Cell(18, 3).Range.Text = "9+15="
$t.Cell(18, 4).Range.Text = "48+4="
$t.Cell(18, 5).Range.Text = "42+19="

$t.Cell(19, 1).Range.Text = "67-8="
$t.Cell(19, 2).Range.Text = "73+19="
$t.Cell(19, 3).Range.Text = "54-45="
$t.Cell(19, 4).Range.Text = "5+79="
$t.Cell(19, 5).Range.Text = "34-6="

$t.Cell(20, 1).Range.Text = "40-19="
$t.Cell(20, 2).Range.Text = "9+19="
$t.Cell(20, 3).Range.Text = "18+47="
$t.Cell(20, 4).Range.Text = "46+37="
$t.Cell(20, 5).Range.Text = "70-11="
